$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest log entry as a plain inline string (matches the
# format used by the other non-header rows in column A).
$ws.Range("A26").Value = "2025-11-15 14:07:47"
